$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing analysis placeholder cells (row 5-8, columns E:H) ---
# Row 5
$ws.Range("E5").Value = "Yes"
$ws.Range("F5").Value = "yes"
$ws.Range("G5").Value = "yes"
$ws.Range("H5").Value = "no"

# Row 6
$ws.Range("E6").Value = "no"
$ws.Range("F6").Value = "yes"
$ws.Range("G6").Value = "yes"
$ws.Range("H6").Value = "yes (Hacking lab has DMARC)"

# --- Add new explanatory row 10 (authored right after the DMARC note above) ---
$ws.Range("H10").Value = "As long as checks are successful, DMARC entry is not relevant, but it seems to be visible in the mail header (X-Spamd) whether DMARC would have been used."
$ws.Range("F10").Value = ""

$ws.Range("F10").WrapText = $true
$ws.Range("H10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 100.8

# Row 7
$ws.Range("E7").Value = "yes"
$ws.Range("F7").Value = "yes"
$ws.Range("G7").Value = "yes"
$ws.Range("H7").Value = "yes"

# Row 8
$ws.Range("E8").Value = "yes"
$ws.Range("F8").Value = "yes"
$ws.Range("G8").Value = "yes (it seems)"
$ws.Range("H8").Value = "not clear"

$ws.Range("H10").Select() | Out-Null
